# Bitacora historica - add June 5th (2020-06-05, serial 43987) data
# across the out_vars / dates_dx / dates_sx / dates_deaths / control_obs sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) out_vars: append row 6, cloning formatting from row 5 then overwriting
#    values.
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Range("A5:J5").Copy($wsOut.Range("A6:J6"))
$wsOut.Range("A6").Value = 43987
$wsOut.Range("B6").Value = 110026
$wsOut.Range("C6").Value = 166049
$wsOut.Range("D6").Value = 48822
$wsOut.Range("E6").Value = 13170
$wsOut.Range("F6").Value = 34.025593950520786
$wsOut.Range("G6").Value = 37437
$wsOut.Range("H6").Value = 3501
$wsOut.Range("I6").Value = 3670
$wsOut.Range("J6").Value = 324897

Write-Output "out_vars done"

# ---------------------------------------------------------------------------
# 2) dates_dx: row 6 already exists pre-formatted (blank), just populate it.
# ---------------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Range("A6").Value = 43987
$wsDx.Range("B6").Value = 0
$wsDx.Range("C6").Value = 1
$wsDx.Range("D6").Value = 1
$wsDx.Range("E6").Value = 1
$wsDx.Range("F6").Value = 0
$wsDx.Range("G6").Value = 0
$wsDx.Range("H6").Value = 0
$wsDx.Range("I6").Value = 4

Write-Output "dates_dx done"

# ---------------------------------------------------------------------------
# 3) dates_sx: append row 6, cloning formatting from row 5.
# ---------------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Range("A5:L5").Copy($wsSx.Range("A6:L6"))
$wsSx.Range("A6").Value = 43987
$wsSx.Range("B6").Value = 0
$wsSx.Range("C6").Value = 1
$wsSx.Range("D6").Value = 0
$wsSx.Range("E6").Value = 1
$wsSx.Range("F6").Value = 1
$wsSx.Range("G6").Value = 1
$wsSx.Range("H6").Value = 0
$wsSx.Range("I6").Value = 1
$wsSx.Range("J6").Value = 1
$wsSx.Range("K6").Value = 0
$wsSx.Range("L6").Value = 0

Write-Output "dates_sx done"

# ---------------------------------------------------------------------------
# 4) dates_deaths: append row 6, cloning formatting from row 5.
# ---------------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Range("A5:H5").Copy($wsDeaths.Range("A6:H6"))
$wsDeaths.Range("A6").Value = 43987
$wsDeaths.Range("B6").Value = 1
$wsDeaths.Range("C6").Value = 0
$wsDeaths.Range("D6").Value = 1
$wsDeaths.Range("E6").Value = 1
$wsDeaths.Range("F6").Value = 1
$wsDeaths.Range("G6").Value = 2
$wsDeaths.Range("H6").Value = 2

Write-Output "dates_deaths done"

# ---------------------------------------------------------------------------
# 5) control_obs: append a new column F (6/5 report), cloning column E's
#    formatting, then overwrite the copied (stale) values/formula.
# ---------------------------------------------------------------------------
$wsCtrl = $wb.Worksheets.Item("control_obs")
$wsCtrl.Range("E1:E20").Copy($wsCtrl.Range("F1:F20"))

$wsCtrl.Range("F1").Value = 43987
$wsCtrl.Range("F2").Value = 3177
$wsCtrl.Range("F3").Value = 2997
$wsCtrl.Range("F4").Value = 2997
$wsCtrl.Range("F5").Value = 2997
$wsCtrl.Range("F6").Value = 2997
$wsCtrl.Range("F7").Value = 2267
$wsCtrl.Range("F8").Value = 4814
$wsCtrl.Range("F9").ClearContents()
$wsCtrl.Range("F10").Value = 145
$wsCtrl.Range("F11").Value = 145
$wsCtrl.Range("F12").Value = 145
$wsCtrl.Range("F13").Value = 145
$wsCtrl.Range("F14").Value = 145
$wsCtrl.Range("F15").Value = 122
$wsCtrl.Range("F16").Value = 157
$wsCtrl.Range("F18").Value = 771
$wsCtrl.Range("F20").Formula = "=SUM(F2:F18)"

Write-Output "control_obs done"
